$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 126.916664
$ws.Range("I53").Value = 117.5
$ws.Range("J53").Value = 136.33333
$ws.Range("K53").Value = 117.5
$ws.Range("L53").Value = 136.33333
$ws.Range("M53").Value = 519.5
$ws.Range("N53").Value = -1410.33333
$ws.Range("H80").Value = 1352.3143
$ws.Range("I80").Value = 1336.5264
$ws.Range("J80").Value = 1371.0625
$ws.Range("K80").Value = 4009.5792
$ws.Range("L80").Value = 4113.1875
$ws.Range("M80").Value = -3011.5792
$ws.Range("N80").Value = -6109.1875
$ws.Range("H83").Value = 1352.3143
$ws.Range("I83").Value = 1336.5264
$ws.Range("J83").Value = 1371.0625
$ws.Range("K83").Value = 12028.7376
$ws.Range("L83").Value = 12339.5625
$ws.Range("M83").Value = -7036.7376
$ws.Range("N83").Value = -22323.5625
$ws.Range("H99").Value = 2811.5
$ws.Range("I99").Value = 1042.8334
$ws.Range("J99").Value = 4580.1665
$ws.Range("K99").Value = 3128.5002
$ws.Range("L99").Value = 13740.4995
$ws.Range("M99").Value = -1630.5002
$ws.Range("N99").Value = -16736.4995
$ws.Range("H100").Value = 1371.6774
$ws.Range("I100").Value = 1040.7142
$ws.Range("J100").Value = 1644.2354
$ws.Range("K100").Value = 1040.7142
$ws.Range("L100").Value = 1644.2354
$ws.Range("M100").Value = -499.7141999999999
$ws.Range("N100").Value = -2726.2354
$ws.Range("H125").Value = 1493
$ws.Range("I125").Value = 792.2
$ws.Range("J125").Value = 4997
$ws.Range("K125").Value = 7129.8
$ws.Range("L125").Value = 44973
$ws.Range("M125").Value = -4669.8
$ws.Range("N125").Value = -49893
$ws.Range("H132").Value = 6079.154
$ws.Range("I132").Value = 5069.4165
$ws.Range("J132").Value = 8351.0625
$ws.Range("K132").Value = 15208.2495
$ws.Range("L132").Value = 25053.1875
$ws.Range("M132").Value = -12678.2495
$ws.Range("N132").Value = -30113.1875
$ws.Range("H137").Value = 5468
$ws.Range("I137").Value = 4809.4443
$ws.Range("J137").Value = 6949.75
$ws.Range("K137").Value = 14428.3329
$ws.Range("L137").Value = 20849.25
$ws.Range("M137").Value = -11878.3329
$ws.Range("N137").Value = -25949.25
$ws.Range("H138").Value = 177785680
$ws.Range("I138").Value = 250006850
$ws.Range("J138").Value = 33343332
$ws.Range("K138").Value = 750020550
$ws.Range("L138").Value = 100029996
$ws.Range("M138").Value = -750015410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4285.375
$ws.Range("I32").Value = 3405.25
$ws.Range("J32").Value = 17487.25
$ws.Range("K32").Value = 3405.25
$ws.Range("L32").Value = 17487.25
$ws.Range("M32").Value = -3118.25
$ws.Range("H88").Value = 1068.1666
$ws.Range("I88").Value = 516.8
$ws.Range("J88").Value = 1280.2307
$ws.Range("K88").Value = 516.8
$ws.Range("L88").Value = 1280.2307
$ws.Range("M88").Value = -110.8
$ws.Range("N88").Value = -2092.2307
$ws.Range("H91").Value = 1068.1666
$ws.Range("I91").Value = 516.8
$ws.Range("J91").Value = 1280.2307
$ws.Range("K91").Value = 516.8
$ws.Range("L91").Value = 1280.2307
$ws.Range("M91").Value = 887.2
$ws.Range("N91").Value = -4088.2307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2777.125
$ws.Range("I99").Value = 2777.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2777.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1279.125
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 4688.4
$ws.Range("I107").Value = 4001.1707
$ws.Range("J107").Value = 7819.1113
$ws.Range("K107").Value = 4001.1707
$ws.Range("L107").Value = 7819.1113
$ws.Range("M107").Value = -2081.1707
$ws.Range("N107").Value = -11659.1113
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3230.6743
$ws.Range("I31").Value = 2331.348
$ws.Range("J31").Value = 4264.9
$ws.Range("K31").Value = 2331.348
$ws.Range("L31").Value = 4264.9
$ws.Range("M31").Value = -2036.348
$ws.Range("N31").Value = -4854.9
$ws.Range("H34").Value = 3230.6743
$ws.Range("I34").Value = 2331.348
$ws.Range("J34").Value = 4264.9
$ws.Range("K34").Value = 2331.348
$ws.Range("L34").Value = 4264.9
$ws.Range("M34").Value = -2129.348
$ws.Range("N34").Value = -4668.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 630.9167
$ws.Range("I5").Value = 662.55554
$ws.Range("J5").Value = 536
$ws.Range("K5").Value = 1987.66662
$ws.Range("L5").Value = 1608
$ws.Range("M5").Value = -1875.66662
$ws.Range("N5").Value = -1832
$ws.Range("H17").Value = 161.875
$ws.Range("I17").Value = 113.57143
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 340.71429
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -171.71429
$ws.Range("N17").Value = -1838
$ws.Range("H39").Value = 12006.538
$ws.Range("I39").Value = 3775
$ws.Range("J39").Value = 13503.182
$ws.Range("K39").Value = 11325
$ws.Range("L39").Value = 40509.546
$ws.Range("M39").Value = -11031
$ws.Range("N39").Value = -41097.546
$ws.Range("H68").Value = 5405.0586
$ws.Range("I68").Value = 2299.7144
$ws.Range("J68").Value = 6210.148
$ws.Range("K68").Value = 6899.1432
$ws.Range("L68").Value = 18630.444
$ws.Range("M68").Value = -6088.1432
$ws.Range("N68").Value = -20252.444
$ws.Range("H71").Value = 5405.0586
$ws.Range("I71").Value = 2299.7144
$ws.Range("J71").Value = 6210.148
$ws.Range("K71").Value = 20697.4296
$ws.Range("L71").Value = 55891.332
$ws.Range("M71").Value = -16641.4296
$ws.Range("N71").Value = -64003.332
$ws.Range("H97").Value = 591.36365
$ws.Range("I97").Value = 543
$ws.Range("J97").Value = 649.4
$ws.Range("K97").Value = 1629
$ws.Range("L97").Value = 1948.2
$ws.Range("M97").Value = -1133
$ws.Range("N97").Value = -2940.2
$ws.Range("H117").Value = 954.4286
$ws.Range("I117").Value = 1115.3334
$ws.Range("J117").Value = 833.75
$ws.Range("K117").Value = 3346.0002
$ws.Range("L117").Value = 2501.25
$ws.Range("M117").Value = 95.99980000000005
$ws.Range("N117").Value = -9385.25
$ws.Range("H131").Value = 7237538
$ws.Range("I131").Value = 15041255
$ws.Range("J131").Value = 4636299
$ws.Range("K131").Value = 45123765
$ws.Range("L131").Value = 13908897
$ws.Range("M131").Value = -45118725
$ws.Range("N131").Value = -13918977
$ws.Range("H135").Value = 630.9167
$ws.Range("I135").Value = 662.55554
$ws.Range("J135").Value = 536
$ws.Range("K135").Value = 5962.99986
$ws.Range("L135").Value = 4824
$ws.Range("M135").Value = -3427.99986
$ws.Range("N135").Value = -9894
$ws.Range("H137").Value = 6758.625
$ws.Range("I137").Value = 5388
$ws.Range("J137").Value = 8129.25
$ws.Range("K137").Value = 16164
$ws.Range("L137").Value = 24387.75
$ws.Range("M137").Value = -11064
$ws.Range("N137").Value = -34587.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 796.6111
$ws.Range("I97").Value = 646.25
$ws.Range("J97").Value = 1097.3334
$ws.Range("K97").Value = 646.25
$ws.Range("L97").Value = 1097.3334
$ws.Range("M97").Value = -150.25
$ws.Range("N97").Value = -2089.3334
$ws.Range("H132").Value = 4670.136
$ws.Range("I132").Value = 4314.8823
$ws.Range("J132").Value = 5878
$ws.Range("K132").Value = 12944.6469
$ws.Range("L132").Value = 17634
$ws.Range("M132").Value = -10414.6469
$ws.Range("N132").Value = -22694

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6415.5
$ws.Range("I40").Value = 4999.2856
$ws.Range("J40").Value = 11372.25
$ws.Range("K40").Value = 4999.2856
$ws.Range("L40").Value = 11372.25
$ws.Range("M40").Value = -4863.2856
$ws.Range("N40").Value = -11644.25
$ws.Range("H58").Value = 23170.727
$ws.Range("I58").Value = 21659.75
$ws.Range("J58").Value = 27200
$ws.Range("K58").Value = 21659.75
$ws.Range("L58").Value = 27200
$ws.Range("M58").Value = -21399.75
$ws.Range("N58").Value = -27720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5012.375
$ws.Range("I113").Value = 4719.8
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 14159.4
$ws.Range("L113").Value = 16500
$ws.Range("M113").Value = -11989.4
$ws.Range("N113").Value = -20840
$ws.Range("H132").Value = 2550
$ws.Range("I132").Value = 1100
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 3300
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -770
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 43208.85
$ws.Range("I136").Value = 55422.848
$ws.Range("J136").Value = 20525.715
$ws.Range("K136").Value = 166268.544
$ws.Range("L136").Value = 61577.145
$ws.Range("M136").Value = -163718.544
$ws.Range("N136").Value = -66677.145
